$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write a "plain" (unstyled / default-style) inline string value into a
# currently-empty cell so that it ends up using the default cell style (index 0),
# matching cells that have no explicit style ("s") attribute.
function Set-PlainCellValue {
    param($Worksheet, $Address, $Text)

    $range = $Worksheet.Range($Address)
    $range.WrapText = $false
    $range.HorizontalAlignment = 1       # xlGeneral
    $range.VerticalAlignment = -4107     # xlBottom
    $range.Value = $Text
}

# --- Row 55: S07_G02_TB001 ---------------------------------------------------
Set-PlainCellValue $ws "F55" "Positions are fetched from Zerodha via ZerodhaClient.list_positions and cached in the positions table; holdings are currently fetched on-demand without DB caching."
$ws.Range("G55").Value = "implemented"
Set-PlainCellValue $ws "H55" "Positions caching covers the main use case; holdings can be added to the cache later if needed."
Set-PlainCellValue $ws "I55" "Add DB caching for holdings if long-term historical holdings views or offline analytics require it."

# --- Row 56: S07_G02_TB002 ---------------------------------------------------
Set-PlainCellValue $ws "F56" "Exposed REST APIs under /api/positions for listing cached positions, triggering a sync, and fetching live holdings."
$ws.Range("G56").Value = "implemented"
Set-PlainCellValue $ws "H56" "Frontend services consume these endpoints for the Positions and Holdings pages."
Set-PlainCellValue $ws "I56" "Extend APIs with filters (e.g., by symbol or product) as usage patterns emerge."

# --- Row 57: S07_G02_TF003 ---------------------------------------------------
Set-PlainCellValue $ws "F57" "Added dedicated Positions and Holdings pages with basic P&L calculations and a manual Refresh for positions."
$ws.Range("G57").Value = "implemented"
Set-PlainCellValue $ws "H57" "User can now inspect positions (with cached P&L) and holdings (with unrealized P&L) directly in the UI."
Set-PlainCellValue $ws "I57" "Enhance UI with aggregate metrics (e.g., total P&L) and filtering in future analytics-focused sprints."
